$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows before the existing row 318 (shifts old 318..336 down to 325..343)
$ws.Rows("318:324").Insert()

# Constant columns shared by every data row in this sheet
$ws.Range("A318:A324").Value2 = 10
$ws.Range("B318:B324").Value2 = "Vega Modelo de Temuco"
$ws.Range("C318:C324").Value2 = "La Araucanía"
$ws.Range("E318:E324").Value2 = 9
$ws.Range("F318:F324").Value2 = "Fruta"
$ws.Range("G318:G324").Value2 = 100103
$ws.Range("H318:H324").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I318:I324").Value2 = 100103001
$ws.Range("J318:J324").Value2 = "Cereza"

# Date column (week of 2023-01-05 == serial 44931), keeps the inherited date style
$ws.Range("D318:D324").Value2 = 44931

# Row 318: Bing / Primera
$ws.Range("K318").Value2 = "Bing"
$ws.Range("L318").Value2 = "Primera"
$ws.Range("M318").Value2 = 125
$ws.Range("N318").Value2 = 500
$ws.Range("O318").Value2 = 500
$ws.Range("P318").Value2 = 500
$ws.Range("Q318").Value2 = "$/kilo (en caja de 15 kilos)"
$ws.Range("R318").Value2 = "Región del Maule"
$ws.Range("S318").Value2 = 500
$ws.Range("T318").Value2 = 1

# Row 319: Bing / Segunda
$ws.Range("K319").Value2 = "Bing"
$ws.Range("L319").Value2 = "Segunda"
$ws.Range("M319").Value2 = 180
$ws.Range("N319").Value2 = 300
$ws.Range("O319").Value2 = 300
$ws.Range("P319").Value2 = 300
$ws.Range("Q319").Value2 = "$/kilo (en caja de 15 kilos)"
$ws.Range("R319").Value2 = "Región del Maule"
$ws.Range("S319").Value2 = 300
$ws.Range("T319").Value2 = 1

# Row 320: Bing / Tercera
$ws.Range("K320").Value2 = "Bing"
$ws.Range("L320").Value2 = "Tercera"
$ws.Range("M320").Value2 = 155
$ws.Range("N320").Value2 = 200
$ws.Range("O320").Value2 = 200
$ws.Range("P320").Value2 = 200
$ws.Range("Q320").Value2 = "$/kilo (en caja de 15 kilos)"
$ws.Range("R320").Value2 = "Región del Maule"
$ws.Range("S320").Value2 = 200
$ws.Range("T320").Value2 = 1

# Row 321: Lapins / Primera
$ws.Range("K321").Value2 = "Lapins"
$ws.Range("L321").Value2 = "Primera"
$ws.Range("M321").Value2 = 125
$ws.Range("N321").Value2 = 500
$ws.Range("O321").Value2 = 500
$ws.Range("P321").Value2 = 500
$ws.Range("Q321").Value2 = "$/kilo (en caja de 15 kilos)"
$ws.Range("R321").Value2 = "Región del Maule"
$ws.Range("S321").Value2 = 500
$ws.Range("T321").Value2 = 1

# Row 322: Lapins / Segunda
$ws.Range("K322").Value2 = "Lapins"
$ws.Range("L322").Value2 = "Segunda"
$ws.Range("M322").Value2 = 155
$ws.Range("N322").Value2 = 300
$ws.Range("O322").Value2 = 300
$ws.Range("P322").Value2 = 300
$ws.Range("Q322").Value2 = "$/kilo (en caja de 15 kilos)"
$ws.Range("R322").Value2 = "Región del Maule"
$ws.Range("S322").Value2 = 300
$ws.Range("T322").Value2 = 1

# Row 323: Lapins / Tercera
$ws.Range("K323").Value2 = "Lapins"
$ws.Range("L323").Value2 = "Tercera"
$ws.Range("M323").Value2 = 125
$ws.Range("N323").Value2 = 200
$ws.Range("O323").Value2 = 200
$ws.Range("P323").Value2 = 200
$ws.Range("Q323").Value2 = "$/kilo (en caja de 15 kilos)"
$ws.Range("R323").Value2 = "Región del Maule"
$ws.Range("S323").Value2 = 200
$ws.Range("T323").Value2 = 1

# Row 324: Sweet Heart / Primera
$ws.Range("K324").Value2 = "Sweet Heart"
$ws.Range("L324").Value2 = "Primera"
$ws.Range("M324").Value2 = 380
$ws.Range("N324").Value2 = 9000
$ws.Range("O324").Value2 = 9000
$ws.Range("P324").Value2 = 9000
$ws.Range("Q324").Value2 = "$/bandeja 10 kilos"
$ws.Range("R324").Value2 = "Región del Maule"
$ws.Range("S324").Value2 = 900
$ws.Range("T324").Value2 = 10
